$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the sheet, reusing the header formatting
# from the neighboring "sum" header cell (G1) so the new header matches
# the existing bold/centered/bordered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New numeric "Save" value for the data row.
$ws.Range("H2").Value = 0
